# Update the "Förändrad" (Changed) date column (C2:C66) from 2023-10-13
# (serial 45212) to 2023-10-22 (serial 45221) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C66").Value = 45221
